$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "2021" column header (R4) - copy formatting from the existing "2020" column (Q4)
$ws.Range("Q4").Copy($ws.Range("R4"))
$ws.Range("R4").Value = 2021

# New data value for 2021 (R5) - copy formatting from the existing "2020" data cell (Q5)
$ws.Range("Q5").Copy($ws.Range("R5"))
$ws.Range("R5").Value = 42.9

# Move the active selection to R9, matching the author's saved view state
$ws.Range("R9").Select()
